# Sample Project / Main.xlsx - "Rules" sheet
# Rule row 11 (the "R40" rule) has its Rule-name cell (B11) changed from
# the text "R40" to the text "1".
#
# The value "1" looks numeric, so it must be forced to remain text
# (matching the shared-string cell the workbook expects) rather than being
# auto-converted to a number by Excel's normal type inference.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("B11")
$cell.NumberFormat = "@"
$cell.Value = "1"
